# ---------------------------------------------------------------------------
# Decision matrix rework: rename the original sheet, fix a couple of input
# numbers on it, then clone it twice to build the "Zeitzone festlegen" and
# "Curl PHP oder Symfony " matrices, and finally add a blank
# "Kontroller Bereiche " sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Fix up the original sheet (becomes "Zeitklasse") -------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Zeitklasse"

# Aufwand/DateTime changed from 6 to 5
$ws1.Range("B2").Value = 5

# Re-colour the Total row: the "B" (DateTime) total becomes the
# green-filled cell, the "C" (Eigene Zeitklasse) total becomes the
# red-filled cell with a plain (non-coloured) font.
$ws1.Range("B8").Interior.Color = 9359529
$ws1.Range("C8").Font.ThemeColor = 1
$ws1.Range("C8").Interior.Color = 255

$null = $ws1.Range("C8").Select()

# --- 2. Clone sheet1 -> "Zeitzone festlegen" --------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Zeitzone festlegen"

# This matrix has 4 criteria (instead of 3), so make room for one extra
# criteria row and one extra formula row.
$ws2.Rows.Item(5).Insert()
$ws2.Rows.Item(9).Insert()

$ws2.Range("B1").Value = "Zentral"
$ws2.Range("C1").Value = "Meheren Orten"

$ws2.Range("A2").Value = "Aufwand"
$ws2.Range("B2").Value = 8
$ws2.Range("C2").Value = 3
$ws2.Range("E2").Value = 1

$ws2.Range("A3").Value = "Vertrautheit mit Methoden"
$ws2.Range("B3").Value = 10
$ws2.Range("C3").Value = 6
$ws2.Range("E3").Value = 2

$ws2.Range("A4").Value = "Fehleranfälligkeit"
$ws2.Range("B4").Value = 8
$ws2.Range("C4").Value = 4
$ws2.Range("E4").Value = 1

$ws2.Range("A5").Value = "Redundanz"
$ws2.Range("B5").Value = 8
$ws2.Range("C5").Value = 6
$ws2.Range("E5").Value = 1

$ws2.Range("A6").Value = "Resultat:"
$ws2.Range("B6").Formula = "=B2*E2"
$ws2.Range("C6").Formula = "=C2*E2"

$ws2.Range("B7").Formula = "=B3*E3"
$ws2.Range("C7").Formula = "=C3*E3"

$ws2.Range("B8").Formula = "=B4*E4"
$ws2.Range("C8").Formula = "=C4*E4"

$ws2.Range("B9").Formula = "=B5*E5"
$ws2.Range("C9").Formula = "=C5*E5"

$ws2.Range("A10").Value = "Total"
$ws2.Range("B10").Formula = "=SUM(B6:B9)"
$ws2.Range("C10").Formula = "=SUM(C6:C9)"

$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(3).AutoFit()

$null = $ws2.Range("E8").Select()

# --- 3. Clone sheet1 -> "Curl PHP oder Symfony " ----------------------------
$ws1.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Curl PHP oder Symfony "

$ws3.Range("B1").Value = "PHP"
$ws3.Range("C1").Value = "Symfony"

$ws3.Range("B2").Value = 5
$ws3.Range("C2").Value = 5
$ws3.Range("E2").Value = 1

$ws3.Range("A3").Value = "Verbreitungsgrad"
$ws3.Range("B3").Value = 10
$ws3.Range("C3").Value = 1
$ws3.Range("E3").Value = 2

$ws3.Range("A4").Value = "Vertrautheit mit Methoden"
$ws3.Range("B4").Value = 5
$ws3.Range("C4").Value = 5
$ws3.Range("E4").Value = 1

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

$null = $ws3.Range("I13").Select()

# --- 4. New blank sheet -> "Kontroller Bereiche " ---------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Kontroller Bereiche "

$null = $ws4.Range("G31").Select()

# --- 5. Restore sheet1 as the active tab ------------------------------------
$ws1.Activate()
